$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "There will be some kind..." text
# (it is the last paragraph in the body, carrying the _GoBack bookmark).
$target = $d.Paragraphs.Item($d.Paragraphs.Count)

# Insert the new paragraphs directly in front of it. We insert OOXML at a
# collapsed range (Start == Start) so Word splices the new paragraphs in
# *before* the existing paragraph, leaving that paragraph (and its
# bookmarkStart/bookmarkEnd) completely untouched.
$insertPoint = $d.Range($target.Range.Start, $target.Range.Start)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:r><w:t xml:space="preserve">There will be </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>some kind of measure</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> that is consistently used across years or geography to measure the rate (%) or absolute numbers of homelessness</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Question</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r><w:t>For people who experience homelessness in the United States, what are the structural or systemic factors that that are associated with changes in homelessness?</w:t></w:r>
</w:p>
<w:p/>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xmlFrag) | Out-Null

# The original paragraph (still holding the old text plus the bookmark) is
# now pushed down below the newly inserted paragraphs. Swap its text for the
# new closing question while leaving the bookmark/run wrapper alone. Scope
# the Find to that specific paragraph's Range so the (identical-looking)
# text in the newly inserted paragraph above is left untouched.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    "There will be some kind of measure that is consistently used across years or geography to measure the rate (%) or absolute numbers of homelessness",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Does the number (or rate) of people who experience homelessness [in X state in the United States] change as [factor] increases or decreases?",
    2) | Out-Null
